$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column F entirely (no longer used)
$ws.Columns.Item(6).Delete()

# Write header row (unchanged text, rewritten for safety)
$ws.Range("B1").Value = "挖掘铲土运输机械产量_同比增长"
$ws.Range("C1").Value = "挖掘铲土运输机械产量_当期值"
$ws.Range("D1").Value = "挖掘铲土运输机械产量_累计值"
$ws.Range("E1").Value = "挖掘铲土运输机械产量_累计增长"

# Write data rows (reordered + new rows appended)
$ws.Range("A2").Value = "2021-10"
$ws.Range("B2").Value = -4.8
$ws.Range("C2").Value = 47605
$ws.Range("D2").Value = 526113
$ws.Range("E2").Value = 15.1

$ws.Range("A3").Value = "2021-11"
$ws.Range("B3").Value = -12.1
$ws.Range("C3").Value = 51217
$ws.Range("D3").Value = 577280
$ws.Range("E3").Value = 13.1

$ws.Range("A4").Value = "2021-12"
$ws.Range("B4").Value = -12.8
$ws.Range("C4").Value = 55463
$ws.Range("D4").Value = 634376
$ws.Range("E4").Value = 9.3

$ws.Range("A5").Value = "2021-02"
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = 128147
$ws.Range("E5").Value = 105.1

$ws.Range("A6").Value = "2021-03"
$ws.Range("B6").Value = 48.4
$ws.Range("C6").Value = 95385
$ws.Range("D6").Value = 227231
$ws.Range("E6").Value = 74.9

$ws.Range("A7").Value = "2021-04"
$ws.Range("B7").Value = 5.1
$ws.Range("C7").Value = 79552
$ws.Range("D7").Value = 310081
$ws.Range("E7").Value = 52.2

$ws.Range("A8").Value = "2021-05"
$ws.Range("B8").Value = -9
$ws.Range("C8").Value = 64162
$ws.Range("D8").Value = 374651
$ws.Range("E8").Value = 36.7

$ws.Range("A9").Value = "2021-06"
$ws.Range("B9").Value = -15
$ws.Range("C9").Value = 51498
$ws.Range("D9").Value = 425461
$ws.Range("E9").Value = 27

$ws.Range("A10").Value = "2021-07"
$ws.Range("B10").Value = -1
$ws.Range("C10").Value = 41186
$ws.Range("D10").Value = 419183
$ws.Range("E10").Value = 25

$ws.Range("A11").Value = "2021-08"
$ws.Range("B11").Value = -12.7
$ws.Range("C11").Value = 35105
$ws.Range("D11").Value = 434434
$ws.Range("E11").Value = 20.2

$ws.Range("A12").Value = "2021-09"
$ws.Range("B12").Value = -4.6
$ws.Range("C12").Value = 46798
$ws.Range("D12").Value = 478096
$ws.Range("E12").Value = 17.5

$ws.Range("A13").Value = "2022-10"
$ws.Range("B13").Value = -4.8
$ws.Range("C13").Value = 46941
$ws.Range("D13").Value = 461272
$ws.Range("E13").Value = -14.5

$ws.Range("A14").Value = "2022-11"
$ws.Range("B14").Value = -6.9
$ws.Range("C14").Value = 49059
$ws.Range("D14").Value = 510014
$ws.Range("E14").Value = -13.8

$ws.Range("A15").Value = "2022-12"
$ws.Range("B15").Value = -24.6
$ws.Range("C15").Value = 42858
$ws.Range("D15").Value = 552604
$ws.Range("E15").Value = -14.8

$ws.Range("A16").Value = "2022-02"
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = ""
$ws.Range("D16").Value = 104211
$ws.Range("E16").Value = -9.8

$ws.Range("A17").Value = "2022-03"
$ws.Range("B17").Value = -22.6
$ws.Range("C17").Value = 59281
$ws.Range("D17").Value = 163555
$ws.Range("E17").Value = -14.9

$ws.Range("A18").Value = "2022-04"
$ws.Range("B18").Value = -48.1
$ws.Range("C18").Value = 34704
$ws.Range("D18").Value = 197180
$ws.Range("E18").Value = -24.9

$ws.Range("A19").Value = "2022-05"
$ws.Range("B19").Value = -25.9
$ws.Range("C19").Value = 41076
$ws.Range("D19").Value = 247946
$ws.Range("E19").Value = -22.8

$ws.Range("A20").Value = "2022-06"
$ws.Range("B20").Value = -2
$ws.Range("C20").Value = 44426
$ws.Range("D20").Value = 293508
$ws.Range("E20").Value = -20.2

$ws.Range("A21").Value = "2022-07"
$ws.Range("B21").Value = -3.4
$ws.Range("C21").Value = 38203
$ws.Range("D21").Value = 330230
$ws.Range("E21").Value = -18.3

$ws.Range("A22").Value = "2022-08"
$ws.Range("B22").Value = 4
$ws.Range("C22").Value = 38281
$ws.Range("D22").Value = 370094
$ws.Range("E22").Value = -16

$ws.Range("A23").Value = "2022-09"
$ws.Range("B23").Value = -12.1
$ws.Range("C23").Value = 43287
$ws.Range("D23").Value = 413447
$ws.Range("E23").Value = -15.6

$ws.Range("A24").Value = "2023-02"
$ws.Range("B24").Value = ""
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = 77330
$ws.Range("E24").Value = -26.4

$ws.Range("A25").Value = "2023-03"
$ws.Range("B25").Value = 2.2
$ws.Range("C25").Value = 61496
$ws.Range("D25").Value = 138587
$ws.Range("E25").Value = -16.3

$ws.Range("A26").Value = "2023-04"
$ws.Range("B26").Value = 9.6
$ws.Range("C26").Value = 38761
$ws.Range("D26").Value = 177015
$ws.Range("E26").Value = -11.3

$ws.Range("A27").Value = "2023-05"
$ws.Range("B27").Value = -19.1
$ws.Range("C27").Value = 33704
$ws.Range("D27").Value = 210723
$ws.Range("E27").Value = -12.7

$ws.Range("A28").Value = "2023-06"
$ws.Range("B28").Value = -17.6
$ws.Range("C28").Value = 37004
$ws.Range("D28").Value = 247761
$ws.Range("E28").Value = -13.6

# Ensure new rows (13-28) get the same date-column style (s="1") as existing rows
$ws.Range("A2").Copy()
$ws.Range("A13:A28").PasteSpecial(-4122)
$excel.CutCopyMode = $false
